$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.655.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "'2.492.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'586.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'176.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").Value = "'0.340"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").Value = "'4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'2.950.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'67.573.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "'2.489.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'11.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'7.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "'351.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'4.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'70.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("D24").Value = "'4.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'1.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'9.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'2.622.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'0.993"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "'0.0₃0907"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").Value = "'508.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'7.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").Value = "'163.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("D37").Value = "'18.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'18.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "'145.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").Value = "'3.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").Value = "'0.515"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "'0.0743"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'0.586"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").Value = "'1.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
